# Auto-generated edit script: updates Leve profit/price figures
# across ARM, CRP, CUL, LTW sheets per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24849.299
$ws.Range("I32").Value = 25861.113
$ws.Range("K32").Value = 25861.113
$ws.Range("M32").Value = -25574.113
$ws.Range("H74").Value = 4145.5674
$ws.Range("I74").Value = 1430.069
$ws.Range("J74").Value = 13989.25
$ws.Range("K74").Value = 1430.069
$ws.Range("L74").Value = 13989.25
$ws.Range("M74").Value = -556.069
$ws.Range("N74").Value = -15737.25
$ws.Range("H77").Value = 4145.5674
$ws.Range("I77").Value = 1430.069
$ws.Range("J77").Value = 13989.25
$ws.Range("K77").Value = 7150.344999999999
$ws.Range("L77").Value = 69946.25
$ws.Range("M77").Value = -2782.344999999999
$ws.Range("N77").Value = -78682.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2221599.8
$ws.Range("I58").Value = 3248204.8
$ws.Range("J58").Value = 10450.308
$ws.Range("K58").Value = 3248204.8
$ws.Range("L58").Value = 10450.308
$ws.Range("M58").Value = -3248001.8
$ws.Range("N58").Value = -10856.308
$ws.Range("H136").Value = 2221599.8
$ws.Range("I136").Value = 3248204.8
$ws.Range("J136").Value = 10450.308
$ws.Range("K136").Value = 9744614.399999999
$ws.Range("L136").Value = 31350.924
$ws.Range("M136").Value = -9742064.399999999
$ws.Range("N136").Value = -36450.924

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 4000
$ws.Range("J100").Value = 4000
$ws.Range("L100").Value = 12000
$ws.Range("N100").Value = -13622
$ws.Range("H103").Value = 433.66666
$ws.Range("I103").Value = 400.5
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 1201.5
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = -322.5
$ws.Range("N103").Value = -3258
$ws.Range("H109").Value = 1339.6364
$ws.Range("J109").Value = 3286.6667
$ws.Range("L109").Value = 9860.000100000001
$ws.Range("N109").Value = -11940.0001
$ws.Range("H112").Value = 2806.75
$ws.Range("I112").Value = 1613.5
$ws.Range("J112").Value = 4000
$ws.Range("K112").Value = 4840.5
$ws.Range("L112").Value = 12000
$ws.Range("M112").Value = -3732.5
$ws.Range("N112").Value = -14216
$ws.Range("H114").Value = 1149.5834
$ws.Range("I114").Value = 424.375
$ws.Range("J114").Value = 2600
$ws.Range("K114").Value = 1273.125
$ws.Range("L114").Value = 7800
$ws.Range("M114").Value = 1980.875
$ws.Range("N114").Value = -14308
$ws.Range("H115").Value = 1966.4
$ws.Range("I115").Value = 1150.5
$ws.Range("J115").Value = 2510.3333
$ws.Range("K115").Value = 3451.5
$ws.Range("L115").Value = 7530.999899999999
$ws.Range("M115").Value = -2276.5
$ws.Range("N115").Value = -9880.999899999999
$ws.Range("H117").Value = 3300
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 3300
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 9900
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -16784
$ws.Range("H118").Value = 5609
$ws.Range("I118").Value = 2508.5
$ws.Range("J118").Value = 18011
$ws.Range("K118").Value = 7525.5
$ws.Range("L118").Value = 54033
$ws.Range("M118").Value = -6282.5
$ws.Range("N118").Value = -56519
$ws.Range("H120").Value = 30010.666
$ws.Range("I120").Value = 20000
$ws.Range("J120").Value = 35016
$ws.Range("K120").Value = 60000
$ws.Range("L120").Value = 105048
$ws.Range("M120").Value = -55162
$ws.Range("N120").Value = -114724
$ws.Range("H121").Value = 717.0526
$ws.Range("I121").Value = 403.5
$ws.Range("J121").Value = 861.7692
$ws.Range("K121").Value = 1210.5
$ws.Range("L121").Value = 2585.3076
$ws.Range("M121").Value = 99.5
$ws.Range("N121").Value = -5205.3076
$ws.Range("H122").Value = 1075.1538
$ws.Range("I122").Value = 674
$ws.Range("J122").Value = 1170.6666
$ws.Range("K122").Value = 6066
$ws.Range("L122").Value = 10535.9994
$ws.Range("M122").Value = -3616
$ws.Range("N122").Value = -15435.9994
$ws.Range("H132").Value = 2959
$ws.Range("I132").Value = 5111
$ws.Range("J132").Value = 2372.0908
$ws.Range("K132").Value = 45999
$ws.Range("L132").Value = 21348.8172
$ws.Range("M132").Value = -43469
$ws.Range("N132").Value = -26408.8172
$ws.Range("H140").Value = 1979.0256
$ws.Range("I140").Value = 1343.9286
$ws.Range("J140").Value = 3595.6365
$ws.Range("K140").Value = 4031.7858
$ws.Range("L140").Value = 10786.9095
$ws.Range("M140").Value = 1148.2142
$ws.Range("N140").Value = -21146.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3130
$ws.Range("I7").Value = 3132.8333
$ws.Range("J7").Value = 3119.8
$ws.Range("K7").Value = 3132.8333
$ws.Range("L7").Value = 3119.8
$ws.Range("M7").Value = -3020.8333
$ws.Range("N7").Value = -3343.8
$ws.Range("H40").Value = 3332.6072
$ws.Range("I40").Value = 3150.4
$ws.Range("J40").Value = 3788.125
$ws.Range("K40").Value = 3150.4
$ws.Range("L40").Value = 3788.125
$ws.Range("M40").Value = -3014.4
$ws.Range("N40").Value = -4060.125
$ws.Range("H61").Value = 38842.734
$ws.Range("I61").Value = 38725.26
$ws.Range("J61").Value = 39228.715
$ws.Range("K61").Value = 38725.26
$ws.Range("L61").Value = 39228.715
$ws.Range("M61").Value = -38523.26
$ws.Range("N61").Value = -39632.715
$ws.Range("H113").Value = 38842.734
$ws.Range("I113").Value = 38725.26
$ws.Range("J113").Value = 39228.715
$ws.Range("K113").Value = 38725.26
$ws.Range("L113").Value = 39228.715
$ws.Range("M113").Value = -36555.26
$ws.Range("N113").Value = -43568.715
$ws.Range("H122").Value = 5606.6304
$ws.Range("I122").Value = 5608.108
$ws.Range("J122").Value = 5600.5557
$ws.Range("K122").Value = 16824.324
$ws.Range("L122").Value = 16801.6671
$ws.Range("M122").Value = -14374.324
$ws.Range("N122").Value = -21701.6671
$ws.Range("H126").Value = 3130
$ws.Range("I126").Value = 3132.8333
$ws.Range("J126").Value = 3119.8
$ws.Range("K126").Value = 9398.499899999999
$ws.Range("L126").Value = 9359.400000000001
$ws.Range("M126").Value = -6928.499899999999
$ws.Range("N126").Value = -14299.4
$ws.Range("H136").Value = 2504.7144
$ws.Range("I136").Value = 1595.2727
$ws.Range("K136").Value = 4785.8181
$ws.Range("M136").Value = -2235.8181
